# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45171 (2023-09-02) to 45172 (2023-09-03).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 260
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45172
}
